$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = [double]"-6.4504025489028581E-3"
$ws.Range("B1").Value = [double]"0.98066620934905835"
$ws.Range("C1").Value = [double]"-1.19863214902941E-2"
$ws.Range("D1").Value = [double]"-1.6207721772448242E-2"
$ws.Range("E1").Value = [double]"-5.8598968420517216E-3"
$ws.Range("F1").Value = [double]"-3.6597770200784002E-2"
$ws.Range("G1").Value = [double]"-1.169790902106695E-2"
$ws.Range("H1").Value = [double]"-1.4110421963665932E-2"
$ws.Range("I1").Value = [double]"-1.0216880691312283E-2"
$ws.Range("J1").Value = [double]"-1.0193501179635595E-2"
$ws.Range("K1").Value = [double]"-1.3228902677560784E-2"
$ws.Range("L1").Value = [double]"0.97930093618373137"
$ws.Range("M1").Value = [double]"-8.6782117927769881E-3"
$ws.Range("N1").Value = [double]"-8.7754474581409215E-3"
$ws.Range("O1").Value = [double]"-1.5883962289262508E-2"
$ws.Range("P1").Value = [double]"-2.1950572130098878E-2"
$ws.Range("Q1").Value = [double]"-1.2828617042535273E-2"
$ws.Range("R1").Value = [double]"-2.8631710755552448E-2"
$ws.Range("S1").Value = [double]"-1.8545741868695939E-2"
$ws.Range("T1").Value = [double]"-3.3333838903449854E-2"
$ws.Range("U1").Value = [double]"0.98950483623283281"
$ws.Range("V1").Value = [double]"0.99122818234694754"
$ws.Range("W1").Value = [double]"-1.2000092562800069E-2"
$ws.Range("X1").Value = [double]"-1.2919420119327365E-2"
$ws.Range("Y1").Value = [double]"-1.2577950587654114E-2"
$ws.Range("Z1").Value = [double]"-6.2034473540487503E-3"
$ws.Range("AA1").Value = [double]"-2.0667544479631698E-2"
$ws.Range("AB1").Value = [double]"-1.5964217160853788E-2"
$ws.Range("AC1").Value = [double]"0.98888832554917339"
$ws.Range("AD1").Value = [double]"-1.0420416214200824E-2"
$ws.Range("AE1").Value = [double]"-1.5299933778063893E-2"
$ws.Range("AF1").Value = [double]"0.97061163176244436"
$ws.Range("AG1").Value = [double]"-1.198312298118725E-2"
$ws.Range("AH1").Value = [double]"-8.6452473506802192E-3"
$ws.Range("AI1").Value = [double]"-1.0536508892805481E-2"
$ws.Range("AJ1").Value = [double]"-6.4448814527411817E-3"
$ws.Range("AK1").Value = [double]"-8.2213312694181077E-3"
$ws.Range("AL1").Value = [double]"-1.1506005474051242E-2"
$ws.Range("AM1").Value = [double]"-1.7806336610932179E-2"
$ws.Range("AN1").Value = [double]"-2.6768957523729466E-2"
$ws.Range("AO1").Value = [double]"-1.1488410332588133E-2"
$ws.Range("AP1").Value = [double]"-5.2115640975039079E-3"
$ws.Range("AQ1").Value = [double]"0.9953821604692259"
$ws.Range("AR1").Value = [double]"-1.5676637159600743E-2"
$ws.Range("AS1").Value = [double]"-2.1072937624592888E-2"
$ws.Range("AT1").Value = [double]"-7.706164770793434E-3"
$ws.Range("AU1").Value = [double]"-2.2808428278693914E-2"
$ws.Range("AV1").Value = [double]"0.99418707731871458"
$ws.Range("AW1").Value = [double]"-2.0065522596864293E-2"
$ws.Range("AX1").Value = [double]"0.98533409386010395"
$ws.Range("AY1").Value = [double]"-1.4972742081440086E-2"
$ws.Range("AZ1").Value = [double]"-1.8903750597335262E-2"
$ws.Range("BA1").Value = [double]"-2.0650695321250266E-2"
$ws.Range("BB1").Value = [double]"-3.141999179928992E-2"
$ws.Range("BC1").Value = [double]"-9.6608410885072778E-3"
$ws.Range("BD1").Value = [double]"0.98507637274401505"
$ws.Range("BE1").Value = [double]"-1.5835564572503225E-2"
$ws.Range("BF1").Value = [double]"-2.8760826597347699E-2"
$ws.Range("BG1").Value = [double]"-1.035577165001391E-2"
$ws.Range("BH1").Value = [double]"-1.191482878189001E-2"
$ws.Range("BI1").Value = [double]"0.99360519910308487"
$ws.Range("BJ1").Value = [double]"-4.8018761704903727E-3"
$ws.Range("BK1").Value = [double]"0.99023135049095701"
$ws.Range("BL1").Value = [double]"-9.3075512785842033E-3"
$ws.Range("BM1").Value = [double]"-7.8217266874536801E-3"
$ws.Range("BN1").Value = [double]"-7.972971160012143E-3"
$ws.Range("BO1").Value = [double]"-1.3040982048442944E-2"
$ws.Range("BP1").Value = [double]"-1.8312163998991285E-2"
$ws.Range("BQ1").Value = [double]"-2.4533619239460704E-2"
$ws.Range("BR1").Value = [double]"-3.0561830880395137E-2"
$ws.Range("BS1").Value = [double]"-2.0090279845705428E-2"
$ws.Range("BT1").Value = [double]"0.99086974692173224"
$ws.Range("BU1").Value = [double]"-2.1343920016180434E-2"
$ws.Range("BV1").Value = [double]"-1.0832143285207698E-2"
$ws.Range("BW1").Value = [double]"0.9629902111391434"
$ws.Range("BX1").Value = [double]"-1.8932677333976761E-2"
$ws.Range("BY1").Value = [double]"-7.3995506934350294E-3"
$ws.Range("BZ1").Value = [double]"-1.1144692593576771E-2"
$ws.Range("CA1").Value = [double]"-1.4198038996972693E-2"
$ws.Range("CB1").Value = [double]"-2.246470641695452E-2"
$ws.Range("CC1").Value = [double]"0.99083814535127879"
$ws.Range("CD1").Value = [double]"-1.3818485208436294E-2"
$ws.Range("CE1").Value = [double]"-5.9207271515047868E-3"
$ws.Range("CF1").Value = [double]"0.99306132309776118"
$ws.Range("CG1").Value = [double]"-1.8492100471798951E-2"
$ws.Range("CH1").Value = [double]"-1.4125327532368318E-2"
$ws.Range("CI1").Value = [double]"-1.9590092891750864E-2"
$ws.Range("CJ1").Value = [double]"-1.7508131933725479E-2"
$ws.Range("CK1").Value = [double]"-3.1158406788846425E-2"
$ws.Range("CL1").Value = [double]"-1.12551168066724E-2"
$ws.Range("CM1").Value = [double]"-1.3307974771696252E-2"
$ws.Range("CN1").Value = [double]"-1.5499828232168004E-2"
$ws.Range("CO1").Value = [double]"-1.6053430893945873E-2"
$ws.Range("CP1").Value = [double]"-1.3532107186406739E-2"
$ws.Range("CQ1").Value = [double]"-3.4488373302490459E-2"
$ws.Range("CR1").Value = [double]"-1.2261999536137203E-2"
$ws.Range("CS1").Value = [double]"-2.0832137760713586E-2"
$ws.Range("CT1").Value = [double]"-1.4416657991352047E-2"
$ws.Range("CU1").Value = [double]"-1.907986811211293E-2"
$ws.Range("CV1").Value = [double]"-7.4417104277832487E-3"
$ws.Range("A2").Value = [double]"0.99834782686856316"
$ws.Range("B2").Value = [double]"-4.8870282131557316E-3"
$ws.Range("C2").Value = [double]"-2.9207890396417062E-3"
$ws.Range("D2").Value = [double]"0.99598204586807459"
$ws.Range("E2").Value = [double]"0.99859667699210908"
$ws.Range("F2").Value = [double]"-9.1852375283296098E-3"
$ws.Range("G2").Value = [double]"-2.8610648200679202E-3"
$ws.Range("H2").Value = [double]"-3.5208111354968438E-3"
$ws.Range("I2").Value = [double]"-2.5575542585713181E-3"
$ws.Range("J2").Value = [double]"-2.6188735429413319E-3"
$ws.Range("K2").Value = [double]"-3.268914715008224E-3"
$ws.Range("L2").Value = [double]"-5.1787897668174676E-3"
$ws.Range("M2").Value = [double]"0.99779193709872238"
$ws.Range("N2").Value = [double]"-2.2189817114851907E-3"
$ws.Range("O2").Value = [double]"-4.01519515537542E-3"
$ws.Range("P2").Value = [double]"-5.5199929126387281E-3"
$ws.Range("Q2").Value = [double]"-3.2378934701796642E-3"
$ws.Range("R2").Value = [double]"-7.2107141158262679E-3"
$ws.Range("S2").Value = [double]"-4.5881732976117895E-3"
$ws.Range("T2").Value = [double]"0.99162969740169227"
$ws.Range("U2").Value = [double]"-2.6476250718255917E-3"
$ws.Range("V2").Value = [double]"-2.1071514086722777E-3"
$ws.Range("W2").Value = [double]"-2.9493820864654103E-3"
$ws.Range("X2").Value = [double]"-3.3032505499956994E-3"
$ws.Range("Y2").Value = [double]"-3.1856300211065918E-3"
$ws.Range("Z2").Value = [double]"-1.5149163330558828E-3"
$ws.Range("AA2").Value = [double]"-5.2478276312444005E-3"
$ws.Range("AB2").Value = [double]"-3.9355338786642696E-3"
$ws.Range("AC2").Value = [double]"-2.7694454629834725E-3"
$ws.Range("AD2").Value = [double]"-2.6243314317434471E-3"
$ws.Range("AE2").Value = [double]"-3.7780447277469273E-3"
$ws.Range("AF2").Value = [double]"-7.4283648451144392E-3"
$ws.Range("AG2").Value = [double]"0.99702276492986963"
$ws.Range("AH2").Value = [double]"-2.1993827629737833E-3"
$ws.Range("AI2").Value = [double]"-2.622621468732593E-3"
$ws.Range("AJ2").Value = [double]"0.99842535021164702"
$ws.Range("AK2").Value = [double]"-2.0047454893189352E-3"
$ws.Range("AL2").Value = [double]"-2.8616383436708809E-3"
$ws.Range("AM2").Value = [double]"-4.3893589062938363E-3"
$ws.Range("AN2").Value = [double]"-6.7515803129924298E-3"
$ws.Range("AO2").Value = [double]"-2.836903517528909E-3"
$ws.Range("AP2").Value = [double]"-1.2848017510258215E-3"
$ws.Range("AQ2").Value = [double]"-1.1006046244702659E-3"
$ws.Range("AR2").Value = [double]"-3.9589122193311064E-3"
$ws.Range("AS2").Value = [double]"-5.2480418735854637E-3"
$ws.Range("AT2").Value = [double]"-1.9147354599405068E-3"
$ws.Range("AU2").Value = [double]"-5.7205467291753728E-3"
$ws.Range("AV2").Value = [double]"-1.4613393204559756E-3"
$ws.Range("AW2").Value = [double]"-4.9668252620808356E-3"
$ws.Range("AX2").Value = [double]"-3.6395847089472247E-3"
$ws.Range("AY2").Value = [double]"0.99625748849497808"
$ws.Range("AZ2").Value = [double]"0.99527582971849382"
$ws.Range("BA2").Value = [double]"-5.2185201715164928E-3"
$ws.Range("BB2").Value = [double]"-7.9466762189207979E-3"
$ws.Range("BC2").Value = [double]"-2.454583142046469E-3"
$ws.Range("BD2").Value = [double]"0.99630497108196692"
$ws.Range("BE2").Value = [double]"-3.9765919017942424E-3"
$ws.Range("BF2").Value = [double]"-7.2301352915351167E-3"
$ws.Range("BG2").Value = [double]"-2.5392163681093056E-3"
$ws.Range("BH2").Value = [double]"-3.0033973956990145E-3"
$ws.Range("BI2").Value = [double]"-1.6222055843765119E-3"
$ws.Range("BJ2").Value = [double]"-1.2286440139371399E-3"
$ws.Range("BK2").Value = [double]"0.99755884921995175"
$ws.Range("BL2").Value = [double]"0.99771614932728314"
$ws.Range("BM2").Value = [double]"-1.9659006979968695E-3"
$ws.Range("BN2").Value = [double]"-1.9975414984590829E-3"
$ws.Range("BO2").Value = [double]"-3.2375241230661973E-3"
$ws.Range("BP2").Value = [double]"0.9954117116984752"
$ws.Range("BQ2").Value = [double]"-6.1958206221436156E-3"
$ws.Range("BR2").Value = [double]"-7.736341512280435E-3"
$ws.Range("BS2").Value = [double]"-4.9634503997779671E-3"
$ws.Range("BT2").Value = [double]"-2.259242393405192E-3"
$ws.Range("BU2").Value = [double]"0.99461801954864459"
$ws.Range("BV2").Value = [double]"0.99729269692026767"
$ws.Range("BW2").Value = [double]"-9.2931192347905471E-3"
$ws.Range("BX2").Value = [double]"-4.8623284612962217E-3"
$ws.Range("BY2").Value = [double]"-1.8462693300767392E-3"
$ws.Range("BZ2").Value = [double]"-2.8089857094220558E-3"
$ws.Range("CA2").Value = [double]"0.99649920024418293"
$ws.Range("CB2").Value = [double]"-5.5846697069008144E-3"
$ws.Range("CC2").Value = [double]"0.99767218444365247"
$ws.Range("CD2").Value = [double]"-3.41305094438991E-3"
$ws.Range("CE2").Value = [double]"-1.4593263032712303E-3"
$ws.Range("CF2").Value = [double]"-1.700708860242543E-3"
$ws.Range("CG2").Value = [double]"-4.7266802823130762E-3"
$ws.Range("CH2").Value = [double]"-3.5888113516663266E-3"
$ws.Range("CI2").Value = [double]"0.99496786625744216"
$ws.Range("CJ2").Value = [double]"-4.390788441401899E-3"
$ws.Range("CK2").Value = [double]"-7.7901743729516916E-3"
$ws.Range("CL2").Value = [double]"0.99722962566637152"
$ws.Range("CM2").Value = [double]"-3.4071835867211485E-3"
$ws.Range("CN2").Value = [double]"-3.9797172045670194E-3"
$ws.Range("CO2").Value = [double]"-3.989911973944423E-3"
$ws.Range("CP2").Value = [double]"-3.292458477051593E-3"
$ws.Range("CQ2").Value = [double]"-8.6650597486387145E-3"
$ws.Range("CR2").Value = [double]"-3.0212371232287865E-3"
$ws.Range("CS2").Value = [double]"-5.3161468958278979E-3"
$ws.Range("CT2").Value = [double]"-3.5824140823284244E-3"
$ws.Range("CU2").Value = [double]"-4.741249721471677E-3"
$ws.Range("CV2").Value = [double]"-1.8115421036353292E-3"
